$wb = $excel.ActiveWorkbook

# "展览" sheet: update "想去人数" (F) counts for rows 3-5
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 124
$wsExhibit.Range("F4").Value = 673
$wsExhibit.Range("F5").Value = 58

# "全部类型" sheet: same events appear one row lower (rows 4-6)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 124
$wsAll.Range("F5").Value = 673
$wsAll.Range("F6").Value = 58
